$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: old row 3 content (Remote Full-Stack Developer / Turing) plus new requirement + link
$ws.Range("A2").Value = "Remote Full-Stack Developer"
$ws.Range("B2").Value = "Turing"
$ws.Range("C2").Value = "['Remote', 'Full-time']"
$ws.Range("D2").Value = "Reposted 1 week ago"
$ws.Range("E2").Value = "* Write clean, reusable, and maintainable code.
* Participate in code reviews to ensure high-quality standards.
* Develop scalable, modular web applications with a focus on security and stability.
* Collaborate with teams to build advanced AI-driven solutions.
* Deliver well-structured and documented code."
$ws.Range("F2").Value = "*   Strong JavaScript or TypeScript skills
*   Bachelor’s or Master’s degree in Computer Science, Engineering, or equivalent experience
*   Strong understanding of ES6 and frameworks like Node.js or React
*   Knowledge of front-end, back-end, or full-stack development
*   Interest in building scalable, secure web applications with clean architecture
*   Good spoken and written communication skills in English
*   Familiarity with additional frameworks like Vue.js, Angular, or Nest.js (Nice to have)
*   Understanding of software quality assurance and test planning (Nice to have)"
$ws.Range("G2").Value = "https://www.linkedin.com/jobs/collections/recommended/?currentJobId=4273383213"

# Row 3: Frontend Development Intern / IG Tech
$ws.Range("A3").Value = "Frontend Development Intern"
$ws.Range("B3").Value = "IG Tech"
$ws.Range("C3").Value = "['Remote', 'Internship']"
$ws.Range("D3").Value = "2 days ago"
$ws.Range("E3").Value = "* Transform designs into functional, visually appealing web applications.
* Build responsive layouts.
* Improve user experiences.
* Work with modern frameworks."
$ws.Range("F3").Value = "*   Pursuing or recently completed a degree in Computer Science, IT, or related fields
*   Proficiency in HTML, CSS, and JavaScript
*   Familiarity with frameworks like React, Angular, or Vue
*   Understanding of responsive design principles
*   Basic knowledge of Git/GitHub for version control
*   Creativity, problem-solving, and attention to detail"
$ws.Range("G3").Value = "https://www.linkedin.com/jobs/collections/recommended/?currentJobId=4300631848"

# Row 4: Web Developer (Wordpress Fullstack Developer) / Uplers
$ws.Range("A4").Value = "Web Developer (Wordpress Fullstack Developer)"
$ws.Range("B4").Value = "Uplers"
$ws.Range("C4").Value = "['₹1.2M/yr - ₹1.8M/yr', 'Remote', 'Full-time']"
$ws.Range("D4").Value = "Reposted 2 weeks ago"
$ws.Range("E4").Value = "*   Design and implement modular, scalable front-end structures.
*   Customize and extend CMS themes and components to meet project requirements.
*   Write semantic, DRY, and well-organized HTML/CSS using naming conventions like BEM.
*   Collaborate with designers, strategists, content, and SEO teams to bring digital concepts to life.
*   Champion best practices and introduce efficiencies in development workflows.
*   Deliver flexible and maintainable front-end codebases that support rapid updates and content scalability.
*   Launch custom, fully editable CMS templates that balance design integrity with editorial flexibility.
*   Implement changes based on technical SEO audits, ensuring site structure, performance, and markup align with search optimization best practices."
$ws.Range("F4").Value = "*   3+ years of experience
*   Proficiency in WordPress themes, Custom PHP (WordPress plugins), and WordPress integration
*   Strong command of JavaScript and PHP
*   Experience with workflow automation and other CMS
*   Ability to architect scalable, reusable components and front-end structures
*   Skilled in writing clean, maintainable, semantic, and DRY HTML/CSS using naming conventions (e.g., BEM), adhering to modern web standards
*   Experience customizing and extending CMS themes and components
*   Ability to champion best practices and introduce efficiencies in development workflows
*   Systematic thinking for reusable patterns and scalable architecture
*   Exceptional attention to detail (visual and functional fidelity)
*   Strong technical communication skills (articulating decisions to non-developers)
*   Proven collaboration skills with cross-functional teams (designers, strategists, content, SEO)
*   Self-directed and proactive in problem-solving
*   Code empathy (writing clean, logical, maintainable code)
*   Adaptability to new tools, workflows, and CMS limitations
*   Process-oriented approach to tasks and delivery
*   User-focused mindset (considering end-user experience)
*   SEO awareness (understanding impact on search performance, technical execution)
*   Design sensitivity (executing visual design with precision)
*   Mentorship mindset"
$ws.Range("G4").Value = "https://www.linkedin.com/jobs/collections/recommended/?currentJobId=4275329833"

# Row 5: Front End Development Intern / PeopleOps Cloud
$ws.Range("A5").Value = "Front End Development Intern"
$ws.Range("B5").Value = "PeopleOps Cloud"
$ws.Range("C5").Value = "['Remote', 'Internship']"
$ws.Range("D5").Value = "4 days ago"
$ws.Range("E5").Value = "* Develop intuitive and scalable user interfaces for web applications.
* Contribute to mobile app development using React Native.
* Assist in designing and implementing conversational flows for chatbot solutions.
* Collaborate with team members to deliver end-to-end product features.
* Ensure UI/UX consistency and responsiveness across devices."
$ws.Range("F5").Value = "*   Proficiency in React.js, React Native, JavaScript/TypeScript, HTML, and CSS.
*   Basic understanding of UI/UX design principles.
*   Familiarity with CSS frameworks (e.g., Bootstrap, Tailwind CSS).
*   Knowledge of WhatsApp APIs or chatbot frameworks (Bonus).
*   Strong problem-solving mindset and ability to work in a dynamic environment.
*   Currently pursuing or recently completed a degree in Computer Science, Engineering, or a related field."
$ws.Range("G5").Value = "https://www.linkedin.com/jobs/collections/recommended/?currentJobId=4301397754"
